# Horarios Linea 141 - scrape update (08:10:18 -> 08:27:16)
# Inserts newly-scraped rows into sheet "LP1912" (4 new rows) and sheet
# "6203-6173" (1 new row), and refreshes the "Ultima actualizacion" /
# "Total filas" header cells on every sheet touched by the scraper run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header refresh
$ws1.Range("A2").Value = "Última actualización: 08:27:16"
$ws1.Range("A3").Value = "Total filas: 65"

# Insert 4 new data rows at their chronological (Hora_Llegada) slots.
# Inserting from the lowest target row to the highest keeps each
# absolute row index valid as we go.
$ws1.Rows.Item(47).Insert()
$ws1.Rows.Item(55).Insert()
$ws1.Rows.Item(63).Insert()
$ws1.Rows.Item(70).Insert()

# New row: 08:27:16 scrape, bus arriving 08:52 at 23_HERNANDEZ
$ws1.Cells.Item(47, 1).Value = "08:27:16"
$ws1.Cells.Item(47, 2).Value = "08:52"
$ws1.Cells.Item(47, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(47, 4).Value = 25
$ws1.Cells.Item(47, 5).Value = "LP1912"

# New row: 08:27:16 scrape, bus arriving 09:17 at 27_EL RETIRO
$ws1.Cells.Item(55, 1).Value = "08:27:16"
$ws1.Cells.Item(55, 2).Value = "09:17"
$ws1.Cells.Item(55, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(55, 4).Value = 50
$ws1.Cells.Item(55, 5).Value = "LP1912"

# New row: 08:27:16 scrape, bus arriving 09:39 at 23_HERNANDEZ
$ws1.Cells.Item(63, 1).Value = "08:27:16"
$ws1.Cells.Item(63, 2).Value = "09:39"
$ws1.Cells.Item(63, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(63, 4).Value = 72
$ws1.Cells.Item(63, 5).Value = "LP1912"

# New row: 08:27:16 scrape, bus arriving 10:13 at 17X38_ROMERO
$ws1.Cells.Item(70, 1).Value = "08:27:16"
$ws1.Cells.Item(70, 2).Value = "10:13"
$ws1.Cells.Item(70, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(70, 4).Value = 106
$ws1.Cells.Item(70, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 (only the scrape timestamp advances this run)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:27:16"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 08:27:16"
$ws3.Range("A3").Value = "Total filas: 8"

# New row: 08:27:16 scrape, bus arriving 10:13 at 215C_LA PLATA
$ws3.Cells.Item(13, 1).Value = "08:27:16"
$ws3.Cells.Item(13, 2).Value = "10:13"
$ws3.Cells.Item(13, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(13, 4).Value = 106
$ws3.Cells.Item(13, 5).Value = "L6203"
